# Update recomputed NATMI edge-weight statistics for the Ccl25-Ccr9 sheet
# (ligand/receptor-expressing cell counts changed from 1 to 3, with
# downstream expression/specificity/edge-weight values recalculated).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.485012999999999
$ws.Range("H2").Value = 16.455039
$ws.Range("I2").Value = 0.2716272065325074
$ws.Range("J2").Value = 0.2716272065325074
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.199522666666667
$ws.Range("N2").Value = 3.598568
$ws.Range("O2").Value = 0.2341312013016629
$ws.Range("P2").Value = 0.2341312013016629
$ws.Range("Q2").Value = 6.579397420461332
$ws.Range("R2").Value = 59.21457678415199
$ws.Range("S2").Value = 0.06359640417167084
$ws.Range("T2").Value = 0.06359640417167085

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.485012999999999
$ws.Range("H3").Value = 16.455039
$ws.Range("I3").Value = 0.2716272065325074
$ws.Range("J3").Value = 0.2716272065325074
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.250404
$ws.Range("N3").Value = 6.751212
$ws.Range("O3").Value = 0.4392495503217397
$ws.Range("P3").Value = 0.4392495503217397
$ws.Range("Q3").Value = 12.343495195252
$ws.Range("R3").Value = 111.091456757268
$ws.Range("S3").Value = 0.1193121283245542
$ws.Range("T3").Value = 0.1193121283245542

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.485012999999999
$ws.Range("H4").Value = 16.455039
$ws.Range("I4").Value = 0.2716272065325074
$ws.Range("J4").Value = 0.2716272065325074
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.063309333333333
$ws.Range("N4").Value = 3.189928
$ws.Range("O4").Value = 0.2075441327510863
$ws.Range("P4").Value = 0.2075441327510863
$ws.Range("Q4").Value = 5.832265516354667
$ws.Range("R4").Value = 52.490389647192
$ws.Range("S4").Value = 0.05637463301138944
$ws.Range("T4").Value = 0.05637463301138944

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.485012999999999
$ws.Range("H5").Value = 16.455039
$ws.Range("I5").Value = 0.2716272065325074
$ws.Range("J5").Value = 0.2716272065325074
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.6100566666666666
$ws.Range("N5").Value = 1.83017
$ws.Range("O5").Value = 0.1190751156255111
$ws.Range("P5").Value = 0.1190751156255111
$ws.Range("Q5").Value = 3.346168747403333
$ws.Range("R5").Value = 30.11551872663
$ws.Range("S5").Value = 0.0323440410248929
$ws.Range("T5").Value = 0.03234404102489291

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 6.542812333333333
$ws.Range("H6").Value = 19.628437
$ws.Range("I6").Value = 0.3240112351547335
$ws.Range("J6").Value = 0.3240112351547334
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.199522666666667
$ws.Range("N6").Value = 3.598568
$ws.Range("O6").Value = 0.2341312013016629
$ws.Range("P6").Value = 0.2341312013016629
$ws.Range("Q6").Value = 7.848251697579554
$ws.Range("R6").Value = 70.63426527821599
$ws.Range("S6").Value = 0.07586113972201332
$ws.Range("T6").Value = 0.07586113972201332

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 6.542812333333333
$ws.Range("H7").Value = 19.628437
$ws.Range("I7").Value = 0.3240112351547335
$ws.Range("J7").Value = 0.3240112351547334
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.250404
$ws.Range("N7").Value = 6.751212
$ws.Range("O7").Value = 0.4392495503217397
$ws.Range("P7").Value = 0.4392495503217397
$ws.Range("Q7").Value = 14.72397104618267
$ws.Range("R7").Value = 132.515739415644
$ws.Range("S7").Value = 0.1423217893409081
$ws.Range("T7").Value = 0.1423217893409081

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.542812333333333
$ws.Range("H8").Value = 19.628437
$ws.Range("I8").Value = 0.3240112351547335
$ws.Range("J8").Value = 0.3240112351547334
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.063309333333333
$ws.Range("N8").Value = 3.189928
$ws.Range("O8").Value = 0.2075441327510863
$ws.Range("P8").Value = 0.2075441327510863
$ws.Range("Q8").Value = 6.957033420281778
$ws.Range("R8").Value = 62.61330078253599
$ws.Range("S8").Value = 0.06724663080179744
$ws.Range("T8").Value = 0.06724663080179742

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.542812333333333
$ws.Range("H9").Value = 19.628437
$ws.Range("I9").Value = 0.3240112351547335
$ws.Range("J9").Value = 0.3240112351547334
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.6100566666666666
$ws.Range("N9").Value = 1.83017
$ws.Range("O9").Value = 0.1190751156255111
$ws.Range("P9").Value = 0.1190751156255111
$ws.Range("Q9").Value = 3.991486282698888
$ws.Range("R9").Value = 35.92337654428999
$ws.Range("S9").Value = 0.03858167529001456
$ws.Range("T9").Value = 0.03858167529001456

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.959658333333334
$ws.Range("H10").Value = 8.878975000000001
$ws.Range("I10").Value = 0.1465673327253718
$ws.Range("J10").Value = 0.1465673327253718
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.199522666666667
$ws.Range("N10").Value = 3.598568
$ws.Range("O10").Value = 0.2341312013016629
$ws.Range("P10").Value = 0.2341312013016629
$ws.Range("Q10").Value = 3.550177256422222
$ws.Range("R10").Value = 31.9515953078
$ws.Range("S10").Value = 0.03431598568257183
$ws.Range("T10").Value = 0.03431598568257184

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.959658333333334
$ws.Range("H11").Value = 8.878975000000001
$ws.Range("I11").Value = 0.1465673327253718
$ws.Range("J11").Value = 0.1465673327253718
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.250404
$ws.Range("N11").Value = 6.751212
$ws.Range("O11").Value = 0.4392495503217397
$ws.Range("P11").Value = 0.4392495503217397
$ws.Range("Q11").Value = 6.660426951966667
$ws.Range("R11").Value = 59.9438425677
$ws.Range("S11").Value = 0.06437963499147638
$ws.Range("T11").Value = 0.06437963499147638

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.959658333333334
$ws.Range("H12").Value = 8.878975000000001
$ws.Range("I12").Value = 0.1465673327253718
$ws.Range("J12").Value = 0.1465673327253718
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.063309333333333
$ws.Range("N12").Value = 3.189928
$ws.Range("O12").Value = 0.2075441327510863
$ws.Range("P12").Value = 0.2075441327510863
$ws.Range("Q12").Value = 3.147032329311112
$ws.Range("R12").Value = 28.3232909638
$ws.Range("S12").Value = 0.0304191899601272
$ws.Range("T12").Value = 0.0304191899601272

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.959658333333334
$ws.Range("H13").Value = 8.878975000000001
$ws.Range("I13").Value = 0.1465673327253718
$ws.Range("J13").Value = 0.1465673327253718
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.6100566666666666
$ws.Range("N13").Value = 1.83017
$ws.Range("O13").Value = 0.1190751156255111
$ws.Range("P13").Value = 0.1190751156255111
$ws.Range("Q13").Value = 1.805559297305555
$ws.Range("R13").Value = 16.25003367575
$ws.Range("S13").Value = 0.01745252209119641
$ws.Range("T13").Value = 0.01745252209119641

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 5.205681333333334
$ws.Range("H14").Value = 15.617044
$ws.Range("I14").Value = 0.2577942255873873
$ws.Range("J14").Value = 0.2577942255873873
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.199522666666667
$ws.Range("N14").Value = 3.598568
$ws.Range("O14").Value = 0.2341312013016629
$ws.Range("P14").Value = 0.2341312013016629
$ws.Range("Q14").Value = 6.244332754776888
$ws.Range("R14").Value = 56.198994792992
$ws.Range("S14").Value = 0.06035767172540686
$ws.Range("T14").Value = 0.06035767172540687

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 5.205681333333334
$ws.Range("H15").Value = 15.617044
$ws.Range("I15").Value = 0.2577942255873873
$ws.Range("J15").Value = 0.2577942255873873
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2.250404
$ws.Range("N15").Value = 6.751212
$ws.Range("O15").Value = 0.4392495503217397
$ws.Range("P15").Value = 0.4392495503217397
$ws.Range("Q15").Value = 11.71488609525867
$ws.Range("R15").Value = 105.433974857328
$ws.Range("S15").Value = 0.113235997664801
$ws.Range("T15").Value = 0.113235997664801

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 5.205681333333334
$ws.Range("H16").Value = 15.617044
$ws.Range("I16").Value = 0.2577942255873873
$ws.Range("J16").Value = 0.2577942255873873
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.063309333333333
$ws.Range("N16").Value = 3.189928
$ws.Range("O16").Value = 0.2075441327510863
$ws.Range("P16").Value = 0.2075441327510863
$ws.Range("Q16").Value = 5.535249548092446
$ws.Range("R16").Value = 49.817245932832
$ws.Range("S16").Value = 0.05350367897777219
$ws.Range("T16").Value = 0.05350367897777219

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 5.205681333333334
$ws.Range("H17").Value = 15.617044
$ws.Range("I17").Value = 0.2577942255873873
$ws.Range("J17").Value = 0.2577942255873873
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.6100566666666666
$ws.Range("N17").Value = 1.83017
$ws.Range("O17").Value = 0.1190751156255111
$ws.Range("P17").Value = 0.1190751156255111
$ws.Range("Q17").Value = 3.175760601942222
$ws.Range("R17").Value = 28.58184541748
$ws.Range("S17").Value = 0.03069687721940724
$ws.Range("T17").Value = 0.03069687721940724
